$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update alpha_distance_range Min: 5.6 -> 5.3
$ws.Range("B2").Value = 5.3

# Update beta_distance_range Min: 5.7 -> 5.4
$ws.Range("B3").Value = 5.4

# Remove the theta_threshold_range row (row 5); rows below shift up.
$ws.Range("A5:C5").Delete()

# The former pie_threshold_range row is now row 5; update its values.
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Match the final selection left by the author's edit.
$ws.Range("B3").Select()

# Page setup as saved by the author (A4 portrait print settings).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
